# Update the "想去人数" (interest count) figures on the 展览 and 全部类型
# sheets to reflect the latest scraped totals (gh-pages output regenerated
# at 456a3b4). Column F (index 6) holds this count on both sheets.

$wb = $excel.ActiveWorkbook

# Updates shared by both sheets (same event rows 4, 8, 14, 19, 20).
$commonUpdates = @(
    @{ Row = 4;  Value = 13472 },
    @{ Row = 8;  Value = 1734 },
    @{ Row = 14; Value = 13474 },
    @{ Row = 19; Value = 8019 },
    @{ Row = 20; Value = 248 }
)

# "展览" (exhibitions) sheet: last updated row is 34.
$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $commonUpdates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}
$ws1.Cells.Item(34, 6).Value = 93

# "全部类型" (all types) sheet aggregates 展览 + 演出 + 本地生活, so the
# equivalent last row is shifted down to 36.
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $commonUpdates) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
$ws4.Cells.Item(36, 6).Value = 93
